$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Semester Information" (sheet1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Semester Information")

# Remove the "Second Semester" and "Other Semester" rows, keep just one
# data row under the header.
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(3).Delete()

# New semester name (plain text).
$ws1.Range("A2").Value = "This is a semester name that should get replaced"

# New start date -- stored as literal text "1/1/1111", NOT as a real date.
# Build it as a text formula result first and paste just the value so Excel
# doesn't reinterpret the digits as a date serial.
$ws1.Range("Z1").Formula = "=""1/1/1111"""
$ws1.Range("Z1").Copy()
$ws1.Range("B2").PasteSpecial(-4163)
$ws1.Range("Z1").ClearContents()

# New end date -- a real date value, formatted with the classic date format.
$ws1.Range("C2").Value = [DateTime]"2050-12-12"
$ws1.Range("C2").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------
# Sheet "Hiatt Staff Emails" (sheet2)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hiatt Staff Emails")

# Fill in the previously-empty row 2 with a new staff email.
$ws2.Range("A2").Value = "schwartzd@brandeis.edu"

# ---------------------------------------------------------------------
# Sheet "Appointment Type Summation" (sheet3)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Appointment Type Summation")

# Drop the old "Appt Type Sum" classification values, leaving just the
# header in B1.
$ws3.Range("B2:B9").ClearContents()

# ---------------------------------------------------------------------
# View state: make "Appointment Type Summation" the active sheet/tab,
# with B2 selected, and leave A3 selected on "Semester Information".
# ---------------------------------------------------------------------
$ws1.Range("A3").Select()
$ws3.Activate()
$ws3.Range("B2").Select()
